$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells: I1 = "I0", J1 = "IF"
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header formatting used by the existing header row (H1 etc.)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New data columns I (I0) and J (IF) for rows 2-20
$data = @(
    @(7, 9),
    @(3, 5),
    @(3, 6),
    @(4, 5),
    @(1, 6),
    @(4, 7),
    @(1, 5),
    @(9, 9),
    @(6, 7),
    @(7, 7),
    @(1, 3),
    @(7, 8),
    @(6, 8),
    @(4, 7),
    @(5, 8),
    @(5, 5),
    @(1, 3),
    @(1, 2),
    @(4, 4)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row++
}
